# Course Management sheet rework: add an employee-number column and a
# resident-registration-number column around the existing 이름/비밀번호
# columns, and fill in a sample row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data row (row 2), name typed in first ------------------------------
$ws.Range("B2").Value = "최소언"

# --- Header row (row 1) -------------------------------------------------
$ws.Range("A1").Value = "직원번호"       # new: employee number
$ws.Range("B1").Value = "이름"           # unchanged
$ws.Range("C1").Value = "주민등록번호"   # new: resident registration number
$ws.Range("D1").Value = "비밀번호"       # unchanged (was column B)

# --- Remainder of the data row --------------------------------------------
$ws.Range("A2").Value = 123
$ws.Range("C2").Value = 123456789
$ws.Range("D2").Value = "thdjs07"

# --- Page margins nudged slightly during the edit -----------------------
$ws.PageSetup.LeftMargin = 50.38
$ws.PageSetup.RightMargin = 50.38

# --- Selection left on B3 after the edit ---------------------------------
$ws.Range("B3").Select()
